$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("4:5").Insert(-4121, 0)
